$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 78; existing rows 78.. shift down by 2.
$ws.Rows("78:79").Insert()

# New row 78 data (Inferno / Primera)
$ws.Cells.Item(78, 1).Value = 8
$ws.Cells.Item(78, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(78, 3).Value = "Coquimbo"
$ws.Cells.Item(78, 4).Value = 44510
$ws.Cells.Item(78, 5).Value = 4
$ws.Cells.Item(78, 6).Value = 100112021
$ws.Cells.Item(78, 7).Value = "Ají"
$ws.Cells.Item(78, 8).Value = "Inferno"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 540
$ws.Cells.Item(78, 11).Value = 23000
$ws.Cells.Item(78, 12).Value = 24000
$ws.Cells.Item(78, 13).Value = 23500
$ws.Cells.Item(78, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(78, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(78, 16).Value = 1958
$ws.Cells.Item(78, 17).Value = 12
$ws.Cells.Item(78, 18).Value = "Hortaliza"

# New row 79 data (Inferno / Segunda)
$ws.Cells.Item(79, 1).Value = 8
$ws.Cells.Item(79, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(79, 3).Value = "Coquimbo"
$ws.Cells.Item(79, 4).Value = 44510
$ws.Cells.Item(79, 5).Value = 4
$ws.Cells.Item(79, 6).Value = 100112021
$ws.Cells.Item(79, 7).Value = "Ají"
$ws.Cells.Item(79, 8).Value = "Inferno"
$ws.Cells.Item(79, 9).Value = "Segunda"
$ws.Cells.Item(79, 10).Value = 400
$ws.Cells.Item(79, 11).Value = 14000
$ws.Cells.Item(79, 12).Value = 15000
$ws.Cells.Item(79, 13).Value = 14500
$ws.Cells.Item(79, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 1208
$ws.Cells.Item(79, 17).Value = 12
$ws.Cells.Item(79, 18).Value = "Hortaliza"
